# Add diary entries for 24/05/2018 and 27/05/2018 ("Add add user dialog")
# at the end of the document, before the final section break.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The "_GoBack" bookmark currently sits right after the last paragraph's
# text. It needs to move into the new final "Entrada" paragraph, between
# the "invitados " and "al crear usuario" runs. Remove it now; it will be
# re-created in the right spot once that paragraph's text exists.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Paragraph: 24/05/2018 (Fechas) ---------------------------------
$ip.InsertParagraphAfter()
$p1 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p1.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:t>24/05/2018</w:t></w:r></w:p>')
$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Paragraph: Boton en el menu para anadir un usuario (Entrada) ---
$ip.InsertParagraphAfter()
$p2 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p2.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Boton</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en el menú para añadir un usuario</w:t></w:r></w:p>')
$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Paragraph: Dialogo para introducir el email para enviar (Entrada) ---
$ip.InsertParagraphAfter()
$p3 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p3.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t>Dialogo para introducir el email para enviar</w:t></w:r></w:p>')
$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Paragraph: 27/05/2018 (Fechas) ---------------------------------
$ip.InsertParagraphAfter()
$p4 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p4.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Fechas"/></w:pPr><w:r><w:t>27/05/2018</w:t></w:r></w:p>')
$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Paragraph: Anadir grupos invitados al crear usuario (Entrada) --
$ip.InsertParagraphAfter()
$p5 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p5.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Entrada"/></w:pPr><w:r><w:t xml:space="preserve">Añadir grupos </w:t></w:r><w:r><w:t xml:space="preserve">invitados </w:t></w:r><w:r><w:t>al crear usuario</w:t></w:r></w:p>')

# Re-insert the _GoBack bookmark between "invitados " and "al crear usuario".
$p5para = $d.Paragraphs($d.Paragraphs.Count)
$bmPos = $p5para.Range.Start + [int]"Añadir grupos invitados ".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$ip = $d.Paragraphs($d.Paragraphs.Count).Range
$ip.Collapse(0)

# --- Final empty paragraph (Entrada, numbering explicitly removed) --
$ip.InsertParagraphAfter()
$p6 = $d.Paragraphs($d.Paragraphs.Count).Range
$null = $p6.InsertXML('<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Entrada"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p>')
